# MHD2-259: Report template and related changes for reporting on 136 genes
#
# Recolour the header table at the top of the clinical-context document:
# the table's default shading (used as the row/table background) and the
# title cell's explicit shading both move from their old colours to the
# new unified background colour ECEAF2.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# wdColor values are 0x00BBGGRR (BGR order), so ECEAF2 -> 0xF2EAEC.
$newColor = 0xF2EAEC

# Table-wide default shading (<w:tblPr><w:shd .../>) - was FFF2CC.
$t.Rows.Shading.BackgroundPatternColor = $newColor

# First cell's explicit shading (<w:tcPr><w:shd .../>) - was E8E7EC.
$t.Cell(1, 1).Shading.BackgroundPatternColor = $newColor
